# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  ("Office Theme" colors) - only wired to the Notes Master
#   ppt/theme/theme2.xml  ("Integral" colors)     - wired to the Slide Master (the
#                                                    theme that actually paints every
#                                                    slide/layout in the deck)
#
# The authored commit swaps the two theme bodies: the visible/applied theme
# (theme2.xml, reached from every slide via the Slide Master) goes from the
# green "Integral" palette to the generic "Office Theme" palette, while the
# "Integral" palette ends up parked in the otherwise-unused theme1.xml slot.
#
# PowerPoint's object model only exposes live, persisted edits to the
# *applied* theme's colors through ThemeColorScheme.Colors(i).RGB (VBA/COM
# RGB long -> 0x00BBGGRR). So we repaint the applied theme's 12 theme colors
# to the stock "Office Theme" values, reproducing the user-visible effect of
# the swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Order is fixed by the OOXML clrScheme sequence:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# RGB isn't available as a builtin here, so colors are passed as the OLE
# "COLORREF" long (0x00BBGGRR) that PowerPoint's own RGB longs use.
$tcs.Colors(1).RGB  = 0x000000    # dk1      000000
$tcs.Colors(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444    # dk2      44546A
$tcs.Colors(4).RGB  = 0xE6E6E7    # lt2      E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B    # accent1  5B9BD5
$tcs.Colors(6).RGB  = 0x317DED    # accent2  ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF    # accent4  FFC000
$tcs.Colors(9).RGB  = 0xC47244    # accent5  4472C4
$tcs.Colors(10).RGB = 0x47AD70    # accent6  70AD47
$tcs.Colors(11).RGB = 0xC16305    # hlink    0563C1
$tcs.Colors(12).RGB = 0x724F95    # folHlink 954F72
